$d = $word.ActiveDocument

# Step 1: delete the "They might be doubtful..." paragraph entirely (para index 4)
$pDoubt = $d.Paragraphs.Item(4)
$pDoubt.Range.Delete()

# Step 2: replace paragraph 3 ("(1)" detail paragraph) runs with the new earnings text
$p3 = $d.Paragraphs.Item(3)
$xmlP3 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr><w:t>The earnings</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr><w:t xml:space="preserve"> after taxes</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr><w:t xml:space="preserve"> is given by $6,068 - $2,185 = $3,883 million.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p3.Range.InsertXML($xmlP3)

# Step 3: replace paragraph 5 ("Actual disbursement..." paragraph) with "Amounts are in millions of $"
$p5 = $d.Paragraphs.Item(5)
$xmlAmounts = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr><w:t>Amounts are in millions of $</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p5.Range.InsertXML($xmlAmounts)

# Step 4: insert the table + trailing paragraphs ((3) and explanation) right before
# the paragraph that now immediately follows paragraph 5 (collapsing to the START of
# paragraph 6 avoids swallowing paragraph 5's own content, unlike collapsing paragraph
# 5's range to its end which merges into the following empty paragraph).
$p6 = $d.Paragraphs.Item(6)
$rngStart = $p6.Range
$rngStart.Collapse(1)
$xmlTail = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="4390"/><w:gridCol w:w="2551"/><w:gridCol w:w="2075"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="4390" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr><w:t>Tax Expense</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr><w:t xml:space="preserve"> To Cash</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr><w:t xml:space="preserve"> To Deferred Tax Liability</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2551" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="right"/><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr><w:t>2,185</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2075" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:jc w:val="right"/><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr><w:t>1,950</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="right"/><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr><w:t>235</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr><w:t>(3)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr><w:t xml:space="preserve">This is because tax authorities and accounting authorities differs in the way they account for income. For example, the revenue in tax authority is decided when cash flows in and the revenue in accounting is not assumed until the delivery has been completed. These leads to discrepancy and finally lead to </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Lucida Sans Typewriter" w:hAnsi="Lucida Sans Typewriter"/></w:rPr><w:t>different income shown to them and shareholders.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rngStart.InsertXML($xmlTail)

Write-Output "Edit complete"
